{"js": "// The site footer block (\"Ver no Jupiter ...\", the \"\u00a9 2020 ...\" copyright\n// line) was removed from the bottom of the document, along with the blank\n// paragraph that had separated it from the \"Requisitos\" list above it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer text paragraphs by their content.\nconst jupiterIdx = items.findIndex((p) => p.text.indexOf(\"Ver no Jupiter\") !== -1);\nconst copyrightIdx = items.findIndex((p) => p.text.indexOf(\"Powered by Jekyll\") !== -1);\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// The blank paragraph immediately before the \"Ver no Jupiter ...\" line is\n// part of the removed block (it only separated the footer from the\n// \"Requisitos\" section above).\nconst blankBeforeIdx = jupiterIdx - 1;\n\nconst toDelete = [];\nif (blankBeforeIdx >= 0 && items[blankBeforeIdx].text === \"\") {\n  toDelete.push(items[blankBeforeIdx]);\n}\ntoDelete.push(items[jupiterIdx]);\ntoDelete.push(items[copyrightIdx]);\n\n// Delete from the bottom up so earlier indices stay valid.\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# The site footer block (\"Ver no Jupiter ...\", the \"\u00a9 2020 ...\" copyright\n# line) is removed from the bottom of the document, along with the blank\n# paragraph that had separated it from the \"Requisitos\" list above it.\n$d = $word.ActiveDocument\n\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($t -like \"*Powered by Jekyll*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -eq -1 -or $copyrightIdx -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# The blank paragraph immediately before the \"Ver no Jupiter ...\" line is\n# part of the removed block (it only separated the footer from the\n# \"Requisitos\" section above).\n$blankIdx = $jupiterIdx - 1\n$blankText = $d.Paragraphs.Item($blankIdx).Range.Text.Trim()\n\n$toDelete = @()\nif ($blankIdx -ge 1 -and $blankText -eq \"\") {\n    $toDelete += $blankIdx\n}\n$toDelete += $jupiterIdx\n$toDelete += $copyrightIdx\n\n# Delete from the bottom up so earlier indices stay valid.\n$toDelete = $toDelete | Sort-Object -Descending\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
